# Auto-generated edit script
# Commit: "Updated cryptos list on Tue Dec  5 03:38:54 UTC 2023 with GitHub Actions"
#
# Updates the crypto price/volume (and two coin-name/link cells that swapped
# rank position) table on Sheet1, matching the target OOXML diff exactly.
#
# Cells are forced to Text format before the write and the style is reset to
# "Normal" afterwards so that numeric-looking strings (e.g. "232.32") are
# stored as literal text -- matching the source workbook's inlineStr cells --
# instead of being auto-coerced to floating point numbers by Excel, and so
# that no stray cell style survives the round trip.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue([string]$cellRef, [string]$value) {
    $rng = $ws.Range($cellRef)
    $rng.NumberFormat = '@'
    $rng.Value = $value
    $rng.Style = 'Normal'
}

Set-TextValue 'D2' '41.813.09'
Set-TextValue 'E2' '  +2.36%  '
Set-TextValue 'D3' '2.233.61'
Set-TextValue 'E3' '  +0.83%  '
Set-TextValue 'E4' '  -0.15%  '
Set-TextValue 'D5' '232.32'
Set-TextValue 'E5' '  +1.13%  '
Set-TextValue 'E6' '  -1.37%  '
Set-TextValue 'D7' '60.83'
Set-TextValue 'E7' '  -5.75%  '
Set-TextValue 'E8' '  -0.08%  '
Set-TextValue 'E9' '  +0.65%  '
Set-TextValue 'D10' '58.06'
Set-TextValue 'E10' '  -1.81%  '
Set-TextValue 'D11' '0.0907'
Set-TextValue 'E11' '  +4.47%  '
Set-TextValue 'E12' '  -0.10%  '
Set-TextValue 'D13' '2.566.84'
Set-TextValue 'E13' '  +0.72%  '
Set-TextValue 'D14' '15.75'
Set-TextValue 'E14' '  -1.08%  '
Set-TextValue 'D15' '22.71'
Set-TextValue 'E15' '  +2.07%  '
Set-TextValue 'D16' '0.806'
Set-TextValue 'E16' '  -2.00%  '
Set-TextValue 'D17' '5.63'
Set-TextValue 'E17' '  +0.10%  '
Set-TextValue 'D18' '2.250.73'
Set-TextValue 'E18' '  +1.29%  '
Set-TextValue 'D19' '41.807.91'
Set-TextValue 'E19' '  +2.79%  '
Set-TextValue 'D20' '72.61'
Set-TextValue 'E20' '  -2.10%  '
Set-TextValue 'D21' '0.0₃0907'
Set-TextValue 'E21' '  +0.41%  '
Set-TextValue 'D22' '6.14'
Set-TextValue 'E22' '  -0.48%  '
Set-TextValue 'D23' '248.23'
Set-TextValue 'E23' '  -1.32%  '
Set-TextValue 'E24' '  -0.07%  '
Set-TextValue 'D25' '2.39'
Set-TextValue 'E25' '  +0.17%  '
Set-TextValue 'D26' '2.31'
Set-TextValue 'E26' '  -2.63%  '
Set-TextValue 'D27' '9.70'
Set-TextValue 'E27' '  -0.10%  '
Set-TextValue 'E28' '  +0.16%  '
Set-TextValue 'D29' '169.70'
Set-TextValue 'E29' '  -1.98%  '
Set-TextValue 'D30' '19.97'
Set-TextValue 'E30' '  -2.05%  '
Set-TextValue 'E31' '  -2.35%  '
Set-TextValue 'D32' '2.67'
Set-TextValue 'E32' '  -5.63%  '
Set-TextValue 'E33' '  -1.32%  '
Set-TextValue 'D34' '5.07'
Set-TextValue 'D35' '4.71'
Set-TextValue 'E35' '  +0.81%  '
Set-TextValue 'D36' '0.0655'
Set-TextValue 'E36' '  +3.63%  '
Set-TextValue 'D37' '6.58'
Set-TextValue 'E37' '  -8.10%  '
Set-TextValue 'D38' '3.64'
Set-TextValue 'E38' '  -4.55%  '
Set-TextValue 'D39' '2.40'
Set-TextValue 'E39' '  -2.49%  '
Set-TextValue 'D40' '0.000240'
Set-TextValue 'E40' '  +17.16%  '
Set-TextValue 'E41' '  -0.12%  '
Set-TextValue 'E42' '  +3.54%  '
Set-TextValue 'D43' '8.63'
Set-TextValue 'E43' '  +0.14%  '
Set-TextValue 'D44' '1.23'
Set-TextValue 'E44' '  -0.32%  '
Set-TextValue 'D45' '98.84'
Set-TextValue 'E45' '  -2.55%  '
Set-TextValue 'E46' '  +1.83%  '
Set-TextValue 'B47' 'Maker'
Set-TextValue 'C47' 'https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr'
Set-TextValue 'D47' '1.473.05'
Set-TextValue 'E47' '  -2.58%  '
Set-TextValue 'B48' 'FTXToken'
Set-TextValue 'C48' 'https://coinranking.com/coin/NfeOYfNcl+ftxtoken-ftt'
Set-TextValue 'D48' '4.35'
Set-TextValue 'E48' '  -11.45%  '
Set-TextValue 'D49' '16.70'
Set-TextValue 'E49' '  -3.70%  '
Set-TextValue 'D50' '2.27'
Set-TextValue 'E50' '  +6.85%  '
Set-TextValue 'E51' '  -2.95%  '
